# EnterpriseJavaTimeLog.xlsx update
# - Adds new Thursday-plan related log entries and a TODO note about
#   per-user profile view authentication (DisplayProfile.java feedback).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 6 new blank rows before row 59 (rows 59-64 become new/blank,
#    everything that used to start at row 59 shifts down to row 65).
$ws.Range("A59:D64").EntireRow.Insert()

# 2) Fill in new content -- written in the same order as the original
#    author typed them, so that the shared-string table indices line up
#    exactly with the target workbook (new strings appended 64..71).
$ws.Range("D53").Value2 = "THURS: above +1, 1 pm - x"
$ws.Range("D55").Value2 = "Thursday plan: "
$ws.Range("D56").Value2 = "go through website and clean up any temporary links or anonymous and hard coded stuff"
$ws.Range("D60").Value2 = "incorporate as much feedback as possible so far"
$ws.Range("D59").Value2 = "deploy to aws to make sure hibernate search doesn't break"
$ws.Range("D57").Value2 = "take an aws snapshot"
$ws.Range("D58").Value2 = "check aws cost structure situation"
$ws.Range("D64").Value2 = "TODO-- PW feedback authentication You may find that to get very fine-grained control over specific data when using a servlet that serve a couple different user types (view my profile versus view somebody else's profile), you may want to check the role within the servlet. The other option would be to set up two servlets ViewMyProfile versus ViewOtherProfile. I can think of pros and cons to each approach..."

# 3) Apply styling to the new D cells to match the rest of column D
#    (wrap text, matching the existing "note" style used throughout).
$ws.Range("D53").WrapText = $true
$ws.Range("D55").WrapText = $true
$ws.Range("D56").WrapText = $true
$ws.Range("D57").WrapText = $true
$ws.Range("D58").WrapText = $true
$ws.Range("D59").WrapText = $true
$ws.Range("D60").WrapText = $true
$ws.Range("D64").WrapText = $true

# 4) Row heights: the newly inserted blank rows 61-63 keep the default
#    15pt height; row 64 (TODO note) needs extra height to show its text.
$ws.Rows(64).RowHeight = 60

# 5) Update the selection to match the author's final cursor position.
$ws.Range("F64").Select()
